$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11 (Marking): Right count and Wrong count corrected
$ws.Range("B11").Value = 4
$ws.Range("C11").Value = -2

# Row 12 (Total): Right count, Wrong count, and Max label corrected
$ws.Range("B12").Value = 12
$ws.Range("C12").Value = -2
$ws.Range("E12").Value = "10 / 112"
